# Updated symbol list on Thu Feb 16 10:47:46 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto-exchange-token rows on Sheet1. The source workbook stores these as
# literal text (e.g. "321.16", "7.37%"), so each target cell is forced to
# Text format before the new literal is written — this prevents Excel's
# automatic type inference from turning a numeric-looking string like
# "321.25" into a real number, or "7.49%" into a percentage value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @{
    "D2"  = "321.25"
    "E2"  = "7.49%"
    "D3"  = "48.75"
    "E3"  = "15.36%"
    "D4"  = "5.261"
    "D5"  = "0.08117"
    "E5"  = "7.62%"
    "D6"  = "4.591"
    "E6"  = "5.34%"
    "D7"  = "1.647"
    "E7"  = "2.69%"
    "D8"  = "1.207"
    "E8"  = "28.89%"
    "D9"  = "0.1295"
    "E9"  = "8.86%"
    "D10" = "0.1947"
    "E10" = "6.80%"
    "D11" = "0.09497"
    "E11" = "4.66%"
    "E12" = "11.87%"
    "D13" = "0.1050"
    "E13" = "0.29%"
    "D14" = "0.001329"
    "E14" = "3.81%"
    "D15" = "0.04166"
    "E15" = "2.04%"
    "D16" = "0.005937"
    "E16" = "1.69%"
    "D17" = "3.341"
    "E17" = "0.02%"
    "D18" = "2.430"
    "E18" = "1.94%"
    "D19" = "0.3404"
    "E19" = "2.09%"
    "D20" = "8.071"
    "E20" = "-3.09%"
    "D21" = "0.1369"
    "E21" = "-2.74%"
    "E22" = "0.89%"
    "D23" = "0.001308"
    "E23" = "3.39%"
    "D24" = "0.004255"
    "E24" = "9.11%"
    "D25" = "0.0001349"
    "E25" = "3.81%"
    "D38" = "0.02720"
    "E38" = "12.46%"
    "D39" = "0.05773"
    "E39" = "10.64%"
    "E40" = "-5.65%"
    "D41" = "0.007698"
    "E41" = "-0.16%"
    "D42" = "0.1444"
    "E42" = "8.79%"
    "E43" = "4.21%"
    "E44" = "14.04%"
    "E45" = "6.58%"
    "D46" = "0.00007023"
    "E46" = "12.85%"
    "E47" = "0.03%"
    "D48" = "0.05500"
    "E48" = "58.00%"
    "E49" = "-4.77%"
    "E50" = "0.03%"
    "E51" = "0.03%"
}

foreach ($cellRef in $cellUpdates.Keys) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
}
